# The commit removes every slide from the deck: the <p:sldIdLst> element
# (and all eight slide parts it referenced) is gone from the saved OOXML,
# leaving only the slide master/layouts/theme untouched.
$p = $ppt.ActivePresentation

# Delete slides from the end backwards so indices stay valid as we go.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}
